# Updates the cryptocurrency price/volume data to the latest scrape.
# Generated from the cell-level diff of the Jul 3 2024 GitHub Actions run.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: cell reference -> new text value.
# Values that look like plain numbers (e.g. prices in column D) are
# prefixed with a leading apostrophe so Excel stores them as literal
# text (matching the workbook's existing text-formatted price column)
# instead of silently converting them to numeric values.
$updates = [ordered]@{
    "D2" = '60.848.40'
    "E2" = '  -3.56%  '
    "D3" = '3.357.15'
    "E3" = '  -2.94%  '
    "D4" = '''1.00'
    "E4" = '  +0.03%  '
    "D5" = '''566.26'
    "E5" = '  -2.21%  '
    "D6" = '''148.15'
    "E6" = '  -0.87%  '
    "E7" = '  +0.18%  '
    "E8" = '  +0.22%  '
    "D9" = '''7.94'
    "E9" = '  +1.01%  '
    "E10" = '  -1.87%  '
    "E11" = '  +0.87%  '
    "D12" = '3.933.11'
    "E12" = '  -2.93%  '
    "E13" = '  +1.04%  '
    "D14" = '''27.97'
    "E14" = '  -2.26%  '
    "D15" = '3.360.15'
    "E15" = '  -2.90%  '
    "E16" = '  -2.05%  '
    "D17" = '60.983.88'
    "E17" = '  -3.41%  '
    "D18" = '''6.35'
    "E18" = '  -1.71%  '
    "D19" = '''14.16'
    "E19" = '  -2.19%  '
    "D20" = '''8.84'
    "E20" = '  -3.67%  '
    "D21" = '''373.78'
    "E21" = '  -3.86%  '
    "E22" = '  +0.89%  '
    "D23" = '''0.559'
    "E23" = '  -0.60%  '
    "E24" = '  +0.00%  '
    "D25" = '3.503.01'
    "E25" = '  -2.54%  '
    "D26" = '''0.0000108'
    "E26" = '  -6.05%  '
    "E27" = '  -4.20%  '
    "D28" = '''1.00'
    "E28" = '  +0.15%  '
    "E29" = '  -4.38%  '
    "E30" = '  +0.01%  '
    "E31" = '  -2.09%  '
    "E32" = '  -5.13%  '
    "D33" = '''22.83'
    "E33" = '  -2.40%  '
    "E34" = '  -3.25%  '
    "D35" = '''5.35'
    "E35" = '  +0.34%  '
    "D36" = '''168.58'
    "E36" = '  -1.01%  '
    "E37" = '  -5.92%  '
    "E38" = '  -4.07%  '
    "D39" = '''29.25'
    "E39" = '  -8.65%  '
    "E40" = '  -2.84%  '
    "D41" = '''0.0750'
    "E41" = '  -3.35%  '
    "D42" = '''42.30'
    "E42" = '  -1.25%  '
    "D43" = '''0.759'
    "E43" = '  -4.57%  '
    "E44" = '  -1.94%  '
    "E45" = '  -4.16%  '
    "E46" = '  -6.27%  '
    "D47" = '2.498.15'
    "E47" = '  -3.44%  '
    "B48" = 'Cosmos'
    "C48" = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
    "D48" = '''6.67'
    "E48" = '  -3.33%  '
    "B49" = 'InjectiveProtocol'
    "C49" = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
    "D49" = '''22.56'
    "E49" = '  -0.56%  '
    "E50" = '  +0.08%  '
    "E51" = '  -2.79%  '
}

foreach ($cellRef in $updates.Keys) {
    $ws.Range($cellRef).Value = $updates[$cellRef]
}
